$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '75.686.44'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +9.41%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.693.33'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +11.27%  '

# Row 4
$ws.Range('E4').Value = '  +0.14%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '189.01'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +14.35%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '589.76'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +4.95%  '

# Row 7
$ws.Range('E7').Value = '  -0.05%  '

# Row 8
$ws.Range('E8').Value = '  +5.82%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.197'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +17.65%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '2.690.82'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +11.11%  '

# Row 11
$ws.Range('E11').Value = '  +1.66%  '

# Row 12
$ws.Range('E12').Value = '  +8.06%  '

# Row 13
$ws.Range('E13').Value = '  +2.20%  '

# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.183.20'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +11.29%  '

# Row 15
$ws.Range('B15').Value = 'WrappedBTC'
$ws.Range('C15').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '75.560.67'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +9.49%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000190'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +7.55%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.73'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +12.00%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.688.00'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +10.95%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.38'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +31.99%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.10'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +12.38%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '376.10'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +10.01%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.30'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +17.18%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.08'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +5.69%  '

# Row 24
$ws.Range('E24').Value = '  +4.82%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.999'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.15%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '70.79'
$ws.Range('D26').ClearFormats()

# Row 27
$ws.Range('E27').Value = '  +10.64%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.50'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +12.60%  '

# Row 29
$ws.Range('E29').Value = '  +11.28%  '

# Row 30
$ws.Range('E30').Value = '  +0.32%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0₃0966'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +14.32%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '523.80'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +16.24%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.42'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +14.24%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.86'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +7.00%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.77'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +10.56%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.998'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -0.01%  '

# Row 37
$ws.Range('E37').Value = '  +9.37%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '162.51'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.24%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '19.36'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +6.57%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '19.39'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +1.68%  '

# Row 41
$ws.Range('E41').Value = '  +0.01%  '

# Row 42
$ws.Range('E42').Value = '  +15.59%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '172.36'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +27.62%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.71'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +13.13%  '

# Row 45
$ws.Range('E45').Value = '  +10.48%  '

# Row 46
$ws.Range('E46').Value = '  +11.89%  '

# Row 47
$ws.Range('E47').Value = '  +16.18%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '39.41'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +4.29%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0852'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +18.25%  '

# Row 50
$ws.Range('E50').Value = '  +8.90%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.543'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +11.74%  '
